$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# H2:H9 used to hold raw pulse-width values. The raw values are moved into a
# new helper column J, and H becomes "=J{row}/2" (half of the raw value).
# Styling (the light "s=2" fill used on rows 4/6/8) follows the value into J.
# ---------------------------------------------------------------------------

# Row 2 - plain (non-shared) formula, no special style involved
$ws.Range("J2").Value = 146.6
$ws.Range("H2").Formula = "=J2/2"

# Row 4 - has the light-fill style; move style + value from H4 to J4
$ws.Range("J4").Value = 120
$ws.Range("H4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("H4").ClearFormats()

# Row 5 - no special style
$ws.Range("J5").Value = 146.6

# Row 6 - has the light-fill style; move style + value from H6 to J6
$ws.Range("J6").Value = 120
$ws.Range("H6").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("H6").ClearFormats()

# Row 7 - no special style
$ws.Range("J7").Value = 140

# Row 8 - has the light-fill style; move style + value from H8 to J8
$ws.Range("J8").Value = 146.6
$ws.Range("H8").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("H8").ClearFormats()

# Row 9 - no special style
$ws.Range("J9").Value = 160

# Fill H3:H9 with the shared formula "=J{row}/2" (H3 is brand new; H4..H9
# already received their value moves above, so this just overwrites them
# with the formula, recreating the shared-formula group that Excel produces
# when a formula is entered once and filled down over a range).
$ws.Range("H3:H9").Formula = "=J3/2"

# ---------------------------------------------------------------------------
# New lists with DO (duty-cycle / pulse) data appended below the existing
# tables, in rows 25 and 26.
# ---------------------------------------------------------------------------
$row25 = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
$row26 = @(0, 0, 2, 6, 10, 11, 8, 4, 2, 0)

for ($i = 0; $i -lt $row25.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(25, $col).Value = $row25[$i]
    $ws.Cells.Item(26, $col).Value = $row26[$i]
}

# ---------------------------------------------------------------------------
# View state: zoomed way out, selection parked on M18.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 40
$null = $ws.Range("M18").Select()
